$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "307.70"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-0.14%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "41.01"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "0.76%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.225"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "1.90%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07656"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "0.43%"
$ws.Range("B6").Value = "GateToken"
$ws.Range("C6").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "4.313"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "1.28%"
$ws.Range("B7").Value = "FTXToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.641"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "1.43%"
$ws.Range("B8").Value = "MXToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.9148"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "1.37%"
$ws.Range("B9").Value = "BTSEToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.436"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "0.29%"
$ws.Range("B10").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C10").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1237"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "12.13%"
$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1824"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "3.35%"
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.09202"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "1.27%"
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.04100"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-2.03%"
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.1051"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "0.00%"
$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001261"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "0.38%"
$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.005875"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "1.14%"
$ws.Range("B17").Value = "UpBots"
$ws.Range("C17").Value = "https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.007509"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "2,395.62%"
$ws.Range("B18").Value = "LEO"
$ws.Range("C18").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.346"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-0.22%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.456"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "13.32%"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "2.11%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.2716"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "1.30%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04043"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-0.59%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001263"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "3.31%"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "4.60%"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-2.18%"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02482"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "4.37%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05337"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "3.15%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.007848"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "1.02%"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "1.01%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.006584"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-2.80%"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-1.93%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.007663"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-12.55%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.3349"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "0.65%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006746"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-3.82%"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.14%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.3722"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "1,101.64%"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-26.09%"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.14%"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "0.14%"
